$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Equal Exchange - One World), which shifts row 3
# (Equal Exchange - Espresso, Decaf) up into row 2.
$ws.Rows(2).Delete()
